$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7428028583526611
$ws.Range("B1").Value = 4.568305492401123
$ws.Range("C1").Value = 5.372705459594727
$ws.Range("D1").Value = 1.161709666252136
$ws.Range("E1").Value = 0.6924885511398315
